# Update loading_percent results for Case_0_8 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4.920233846125796
$ws.Range("D2").Value = 4.26030435454938
$ws.Range("E2").Value = 16.5182167249567
$ws.Range("F2").Value = 22.05448126744788
$ws.Range("G2").Value = 3.610796766715304
$ws.Range("K2").Value = 15.88773969315364
$ws.Range("O2").Value = 19.62758929797896

# Row 3
$ws.Range("C3").Value = 4.749570996603186
$ws.Range("D3").Value = 4.232734541106978
$ws.Range("E3").Value = 15.57829436551835
$ws.Range("F3").Value = 22.10850585961802
$ws.Range("G3").Value = 3.613461599718025
$ws.Range("K3").Value = 15.11944908495187
$ws.Range("O3").Value = 19.7443556213385

# Row 4
$ws.Range("C4").Value = 4.643154055282187
$ws.Range("D4").Value = 4.215922356990685
$ws.Range("E4").Value = 14.9764604527888
$ws.Range("F4").Value = 22.15256032837706
$ws.Range("G4").Value = 3.615180304542782
$ws.Range("K4").Value = 14.62730794762479
$ws.Range("O4").Value = 19.8244267812052

# Row 5
$ws.Range("C5").Value = 4.599456353541459
$ws.Range("D5").Value = 4.209104587827865
$ws.Range("E5").Value = 14.72526202066048
$ws.Range("F5").Value = 22.17322544499681
$ws.Range("G5").Value = 3.615901508794893
$ws.Range("K5").Value = 14.42184242018346
$ws.Range("O5").Value = 19.85914331061602

# Row 6
$ws.Range("C6").Value = 4.592182728632412
$ws.Range("D6").Value = 4.207974637994665
$ws.Range("E6").Value = 14.68320049324652
$ws.Range("F6").Value = 22.17681995903875
$ws.Range("G6").Value = 3.616022523807858
$ws.Range("K6").Value = 14.38743530437854
$ws.Range("O6").Value = 19.86503343661537

# Row 7
$ws.Range("C7").Value = 4.642565969463111
$ws.Range("D7").Value = 4.215830269706748
$ws.Range("E7").Value = 14.97309638397218
$ws.Range("F7").Value = 22.15282807583164
$ws.Range("G7").Value = 3.615189946563786
$ws.Range("K7").Value = 14.62455655031473
$ws.Range("O7").Value = 19.82488655662731

# Row 8
$ws.Range("C8").Value = 4.861772167936813
$ws.Range("D8").Value = 4.250776788128404
$ws.Range("E8").Value = 16.19940750658397
$ws.Range("F8").Value = 22.07083630799838
$ws.Range("G8").Value = 3.611698527025648
$ws.Range("K8").Value = 15.62718845590669
$ws.Range("O8").Value = 19.66610025984131

# Row 9
$ws.Range("C9").Value = 5.275611982874131
$ws.Range("D9").Value = 4.320039924007671
$ws.Range("E9").Value = 18.49317080658181
$ws.Range("F9").Value = 21.99737278499233
$ws.Range("G9").Value = 3.605502868884516
$ws.Range("K9").Value = 17.42355055341829
$ws.Range("O9").Value = 19.42203813090024

# Row 10
$ws.Range("C10").Value = 5.566200427766797
$ws.Range("D10").Value = 4.371125826048504
$ws.Range("E10").Value = 20.14271095942959
$ws.Range("F10").Value = 21.99779073258719
$ws.Range("G10").Value = 3.601342884360899
$ws.Range("K10").Value = 18.63130899603606
$ws.Range("O10").Value = 19.28491773831567

# Row 11
$ws.Range("C11").Value = 5.694820208658134
$ws.Range("D11").Value = 4.394356256296538
$ws.Range("E11").Value = 20.8510766058532
$ws.Range("F11").Value = 22.00998434585059
$ws.Range("G11").Value = 3.59953446158963
$ws.Range("K11").Value = 19.15510940501479
$ws.Range("O11").Value = 19.23194022509221

# Row 12
$ws.Range("C12").Value = 5.742965344609335
$ws.Range("D12").Value = 4.403147283938737
$ws.Range("E12").Value = 21.1133119709038
$ws.Range("F12").Value = 22.016339323989
$ws.Range("G12").Value = 3.598861654173872
$ws.Range("K12").Value = 19.34969142436087
$ws.Range("O12").Value = 19.21324980986692

# Row 13
$ws.Range("C13").Value = 5.732622040504245
$ws.Range("D13").Value = 4.401254310557012
$ws.Range("E13").Value = 21.05710151772959
$ws.Range("F13").Value = 22.01489323846619
$ws.Range("G13").Value = 3.59900602260154
$ws.Range("K13").Value = 19.30795357482319
$ws.Range("O13").Value = 19.21721384486646

# Row 14
$ws.Range("C14").Value = 5.698792618237488
$ws.Range("D14").Value = 4.395079642350309
$ws.Range("E14").Value = 20.87277106402972
$ws.Range("F14").Value = 22.01047229625371
$ws.Range("G14").Value = 3.599478869174303
$ws.Range("K14").Value = 19.17119377561158
$ws.Range("O14").Value = 19.23037496667275

# Row 15
$ws.Range("C15").Value = 5.677996842622285
$ws.Range("D15").Value = 4.39129658532149
$ws.Range("E15").Value = 20.7590820381577
$ws.Range("F15").Value = 22.00799090690289
$ws.Range("G15").Value = 3.599770062356459
$ws.Range("K15").Value = 19.08693121841346
$ws.Range("O15").Value = 19.23861564879398

# Row 16
$ws.Range("C16").Value = 5.55771883473164
$ws.Range("D16").Value = 4.369607094591978
$ws.Range("E16").Value = 20.09557506550816
$ws.Range("F16").Value = 21.99723643691106
$ws.Range("G16").Value = 3.601462752633091
$ws.Range("K16").Value = 18.59655420385972
$ws.Range("O16").Value = 19.28857075138881

# Row 17
$ws.Range("C17").Value = 5.482982506472664
$ws.Range("D17").Value = 4.356295781328588
$ws.Range("E17").Value = 19.67779797484908
$ws.Range("F17").Value = 21.993722367219
$ws.Range("G17").Value = 3.602522619669873
$ws.Range("K17").Value = 18.28909636171711
$ws.Range("O17").Value = 19.32163792986944

# Row 18
$ws.Range("C18").Value = 5.439662526093962
$ws.Range("D18").Value = 4.348638909016681
$ws.Range("E18").Value = 19.43354659497954
$ws.Range("F18").Value = 21.99283049279833
$ws.Range("G18").Value = 3.603140135604568
$ws.Range("K18").Value = 18.10984776299485
$ws.Range("O18").Value = 19.34154086485369

# Row 19
$ws.Range("C19").Value = 5.424939348692758
$ws.Range("D19").Value = 4.346046467583076
$ws.Range("E19").Value = 19.35016658249905
$ws.Range("F19").Value = 21.99272208543215
$ws.Range("G19").Value = 3.603350576384241
$ws.Range("K19").Value = 18.04874677232701
$ws.Range("O19").Value = 19.34843081062401

# Row 20
$ws.Range("C20").Value = 5.490973202650547
$ws.Range("D20").Value = 4.357712884223275
$ws.Range("E20").Value = 19.7226803445474
$ws.Range("F20").Value = 21.99397947314809
$ws.Range("G20").Value = 3.602408977048861
$ws.Range("K20").Value = 18.32207556513947
$ws.Range("O20").Value = 19.31802628435752

# Row 21
$ws.Range("C21").Value = 5.708744695110096
$ws.Range("D21").Value = 4.396893486579438
$ws.Range("E21").Value = 20.92707614691482
$ws.Range("F21").Value = 22.01172360187017
$ws.Range("G21").Value = 3.599339657491422
$ws.Range("K21").Value = 19.21146636470136
$ws.Range("O21").Value = 19.22647186481489

# Row 22
$ws.Range("C22").Value = 5.847786648756218
$ws.Range("D22").Value = 4.422464189836281
$ws.Range("E22").Value = 21.67921913625817
$ws.Range("F22").Value = 22.03345121640704
$ws.Range("G22").Value = 3.597403609175363
$ws.Range("K22").Value = 19.77073201405528
$ws.Range("O22").Value = 19.17463574882485

# Row 23
$ws.Range("C23").Value = 5.773891959539124
$ws.Range("D23").Value = 4.408821425128562
$ws.Range("E23").Value = 21.2809772381904
$ws.Range("F23").Value = 22.02092480213907
$ws.Range("G23").Value = 3.598430539927016
$ws.Range("K23").Value = 19.47427915142602
$ws.Range("O23").Value = 19.20156337690024

# Row 24
$ws.Range("C24").Value = 5.487361705043461
$ws.Range("D24").Value = 4.357072224157274
$ws.Range("E24").Value = 19.70240169223034
$ws.Range("F24").Value = 21.99385972155134
$ws.Range("G24").Value = 3.602460329378906
$ws.Range("K24").Value = 18.30717340698233
$ws.Range("O24").Value = 19.31965633169612

# Row 25
$ws.Range("C25").Value = 5.165778704104032
$ws.Range("D25").Value = 4.301249550244286
$ws.Range("E25").Value = 17.848219912773
$ws.Range("F25").Value = 22.00776869087911
$ws.Range("G25").Value = 3.607109771999675
$ws.Range("K25").Value = 16.9567438679337
$ws.Range("O25").Value = 19.48073636245494
